$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testsheet2")

# 1. Insert a new blank column at A; old A:E (FirstName..Password) shift to B:F
$ws.Columns("A").Insert()

# 2. Fill in the new "TC ID/Name" column (A) header + the first test-case row
$ws.Range("A1").Value = "TC ID/Name"
$ws.Range("A2").Value = "testingValidUserSignUp"

# 3. "Phone" (now column D) becomes "Mobile Number" -- same column, new header text
$ws.Range("D1").Value = "Mobile Number"

# 4. Duplicate data row 2 (now has TC ID/Name..Password values) into rows 3 and 4
$ws.Rows("2").Copy()
$ws.Rows("3").Insert()
$ws.Rows("2").Copy()
$ws.Rows("4").Insert()

# 5. Overwrite the TC ID/Name cell for the two new rows with their own test case names
$ws.Range("A3").Value = "testingInvalidUserSignUp_emailAlreadyExists"
$ws.Range("A4").Value = "testingInvalidUserSignUp_emailWrongFormat"

# Match header formatting (fill/border/left align) used by the other header cells
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# 6. Fill in the new "Expected Alert Message" column (G)
$ws.Range("G1").Value = "Expected Alert Message"
$ws.Range("G3").Value = "Email Already Exists."
$ws.Range("G4").Value = "The Email field must contain a valid email address."

$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("G4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 7. Column widths: new columns + the auto-fit "Mobile Number" header
# (values chosen so the pixel-quantized stored width lands as close as
# possible to the saved widths of 38.875 / 14.25 / 43.375)
$ws.Columns("A").ColumnWidth = 38.0833
$ws.Columns("D").ColumnWidth = 13.41667
$ws.Columns("G").ColumnWidth = 42.5833

# 8. Selection, matching the saved state in the workbook
$ws.Range("G1").Select() | Out-Null
